$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 holds the log date as literal text ("2019-06-15" -> "2019-06-16").
# Writing straight into Value would make Excel auto-detect it as a date
# serial, so stage the text in a scratch cell formatted as Text, then
# paste-special (values only) into B1 so the original General style/type
# of the cell is preserved while the content still lands as a string.
$scratch = $ws.Range("D1")
$scratch.NumberFormat = "@"
$scratch.Value = "2019-06-16"
$scratch.Copy()
$ws.Range("B1").PasteSpecial(-4163)
$scratch.Clear()

# Updated per-level log counts for the refreshed stats snapshot.
$ws.Range("B2").Value = 90.0
$ws.Range("B3").Value = 90.0
$ws.Range("B4").Value = 16.0
$ws.Range("B5").Value = 13.0
$ws.Range("B6").Value = 17.0
$ws.Range("B7").Value = 23.0
$ws.Range("B8").Value = 21.0
